$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12, pushing the old row 12 ("IMPALA") down to row 13
[void]$ws.Rows("12:12").Insert()

# Update row 11 values (PPO 1.2, 1.3, 1.4, 3.1)
$ws.Range("B11").Value = 12
$ws.Range("C11").Value = 3.1
$ws.Range("D11").Value = 0.2
$ws.Range("E11").Value = 23.4
$ws.Range("F11").Value = 14360000
$ws.Range("F11").NumberFormat = $ws.Range("F9").NumberFormat

# Fill in the newly inserted row 12
$ws.Range("A12").Value = "PPO 2.2, 3.1"
$ws.Range("G12").Value = "because lowest successrate on 2.2"

# Clear formatting/values accidentally carried over into the blank cells of row 12
$ws.Range("B12:E12").Clear()

# Update sheet view selection to match the saved state
[void]$ws.Range("K16").Select()
